$wb = $excel.ActiveWorkbook

# --- Step1_Data ---
$ws = $wb.Worksheets.Item("Step1_Data")
$ws.Range("AA4").Value = 0.01002285616973044
$ws.Range("AA5").Value = 0.007795058174530614
$ws.Range("AA6").Value = 0.007203306873766927
$ws.Range("AB2").Value = 0.001532337420433319
$ws.Range("AB3").Value = 0.01794929670189378
$ws.Range("AB4").Value = 0.006842668121921212
$ws.Range("AD2").Value = 0.01211244588483981
$ws.Range("AD4").Value = 0.03947350409632015
$ws.Range("AD5").Value = 0.0103159389425093
$ws.Range("AD6").Value = 0.01357861353158792
$ws.Range("AE2").Value = 0.1025786316457903
$ws.Range("AE3").Value = 0.04968913244502298
$ws.Range("AE4").Value = 0.03352112663135218
$ws.Range("AE5").Value = 0.04135721623161707
$ws.Range("AE6").Value = 0.03801787359746096
$ws.Range("AF2").Value = 0.02825362570385796
$ws.Range("AF3").Value = 0.04068509596301505
$ws.Range("AF4").Value = 0.03128763537187704
$ws.Range("AF5").Value = 0.0612322603389075
$ws.Range("AF6").Value = 0.04717991363016376
$ws.Range("AG2").Value = 0.003583926915926361
$ws.Range("AG3").Value = 0.00705459114029016
$ws.Range("AG5").Value = 0.03419821812055224
$ws.Range("AG6").Value = 0.01649550232964162
$ws.Range("AH4").Value = 0.004689519446374135
$ws.Range("AH5").Value = 0.04207322417519559
$ws.Range("AH6").Value = 0.02682315456678177
$ws.Range("AI3").Value = 0.0009876395324247742
$ws.Range("AJ2").Value = 0.0180246493372624
$ws.Range("D5").Value = 0.02018260650069554
$ws.Range("D6").Value = 0.02051331462385681
$ws.Range("E2").Value = 0.2322838788058234
$ws.Range("E4").Value = 0.09920926753084587
$ws.Range("E5").Value = 0.07574603178797346
$ws.Range("E6").Value = 0.006555557844093176
$ws.Range("F3").Value = 0.2029300480133553
$ws.Range("F4").Value = 0.1256738871696395
$ws.Range("F5").Value = 0.07016063378674665
$ws.Range("G2").Value = 0.2565491549287454
$ws.Range("G3").Value = 0.1743099004836132
$ws.Range("G5").Value = 0.02554441335181604
$ws.Range("G6").Value = 0.02409479471203695
$ws.Range("H3").Value = 0.0243839523784519
$ws.Range("H4").Value = 0.02651577652828388
$ws.Range("H5").Value = 0.1963247361980084
$ws.Range("H6").Value = 0.3584431042379857
$ws.Range("I2").Value = 0.02237471229987421
$ws.Range("I4").Value = 0.04736548656558134
$ws.Range("J3").Value = 0.01612001231296669
$ws.Range("J6").Value = 0.006756873532457942
$ws.Range("K4").Value = 0.0003595980011707987
$ws.Range("L2").Value = 0.03777298937769286
$ws.Range("L4").Value = 0.1306194795028028
$ws.Range("L5").Value = 0.04547375936846133
$ws.Range("L6").Value = 0.08110026839183615
$ws.Range("M2").Value = 0.01574234655326109
$ws.Range("M3").Value = 0.2612027650757788
$ws.Range("M4").Value = 0.2212629084967883
$ws.Range("M5").Value = 0.2044081930082545
$ws.Range("M6").Value = 0.155100368069891
$ws.Range("N2").Value = 0.03840349202796589
$ws.Range("N3").Value = 0.008542583017184489
$ws.Range("N5").Value = 0.04335573074610145
$ws.Range("N6").Value = 0.03115842138224091
$ws.Range("O2").Value = 0.1201366128835211
$ws.Range("O3").Value = 0.005949087554763117
$ws.Range("O6").Value = 0.01321286408151596
$ws.Range("P3").Value = 0.04604159994469997
$ws.Range("Q3").Value = 0.002151683372682401
$ws.Range("Q5").Value = 0.01024316162988733
$ws.Range("Q6").Value = 0.01760163875795359
$ws.Range("R2").Value = 0.005206092100606968
$ws.Range("T4").Value = 0.1165363872464782
$ws.Range("T5").Value = 0.05338319683546135
$ws.Range("T6").Value = 0.06897997570232654
$ws.Range("U2").Value = 0.07985503990305431
$ws.Range("U3").Value = 0.07479109547406096
$ws.Range("U4").Value = 0.08772343615357127
$ws.Range("U5").Value = 0.02863603625349264
$ws.Range("U6").Value = 0.03601227134036555
$ws.Range("V2").Value = 0.01720180795267639
$ws.Range("V3").Value = 0.01645585810400555
$ws.Range("V5").Value = 0.02956958454978924
$ws.Range("V6").Value = 0.03117218279403674
$ws.Range("W3").Value = 0.01608538213304919
$ws.Range("X4").Value = 0.01612827924609601
$ws.Range("Y2").Value = 0.00838825625866833
$ws.Range("Y3").Value = 0.03467027635274147
$ws.Range("Z4").Value = 0.002768183721166776

# --- Step2_Sj ---
$ws = $wb.Worksheets.Item("Step2_Sj")
$ws.Range("AA2").Value = 0.83391438309189
$ws.Range("AA3").Value = 0.883634244217353
$ws.Range("AA4").Value = 0.8841855463321553
$ws.Range("AA5").Value = 0.8108231421912185
$ws.Range("AA6").Value = 0.857904942344364
$ws.Range("AB2").Value = 0.8354467205123234
$ws.Range("AB3").Value = 0.9015835409192469
$ws.Range("AB4").Value = 0.8910282144540765
$ws.Range("AB5").Value = 0.8108231421912185
$ws.Range("AB6").Value = 0.857904942344364
$ws.Range("AC2").Value = 0.8354467205123234
$ws.Range("AC3").Value = 0.9015835409192469
$ws.Range("AC4").Value = 0.8910282144540765
$ws.Range("AC5").Value = 0.8108231421912185
$ws.Range("AC6").Value = 0.857904942344364
$ws.Range("AD2").Value = 0.8475591663971632
$ws.Range("AD3").Value = 0.9015835409192469
$ws.Range("AD4").Value = 0.9305017185503967
$ws.Range("AD5").Value = 0.8211390811337278
$ws.Range("AD6").Value = 0.8714835558759519
$ws.Range("AE2").Value = 0.9501377980429535
$ws.Range("AE3").Value = 0.9512726733642698
$ws.Range("AE4").Value = 0.9640228451817489
$ws.Range("AE5").Value = 0.8624962973653448
$ws.Range("AE6").Value = 0.9095014294734129
$ws.Range("AF2").Value = 0.9783914237468114
$ws.Range("AF3").Value = 0.9919577693272849
$ws.Range("AF4").Value = 0.9953104805536259
$ws.Range("AF5").Value = 0.9237285577042523
$ws.Range("AF6").Value = 0.9566813431035767
$ws.Range("AG2").Value = 0.9819753506627378
$ws.Range("AG3").Value = 0.999012360467575
$ws.Range("AG4").Value = 0.9953104805536259
$ws.Range("AG5").Value = 0.9579267758248046
$ws.Range("AG6").Value = 0.9731768454332183
$ws.Range("AH2").Value = 0.9819753506627378
$ws.Range("AH3").Value = 0.999012360467575
$ws.Range("AH4").Value = 1
$ws.Range("AH6").Value = 1
$ws.Range("AI2").Value = 0.9819753506627378
$ws.Range("AI3").Value = 0.9999999999999998
$ws.Range("AI4").Value = 1
$ws.Range("AI6").Value = 1
$ws.Range("AJ3").Value = 0.9999999999999998
$ws.Range("AJ4").Value = 1
$ws.Range("AJ6").Value = 1
$ws.Range("AK3").Value = 0.9999999999999998
$ws.Range("AK4").Value = 1
$ws.Range("AK6").Value = 1
$ws.Range("D5").Value = 0.02018260650069554
$ws.Range("D6").Value = 0.02051331462385681
$ws.Range("E2").Value = 0.2322838788058234
$ws.Range("E4").Value = 0.09920926753084587
$ws.Range("E5").Value = 0.095928638288669
$ws.Range("E6").Value = 0.02706887246794999
$ws.Range("F2").Value = 0.2322838788058234
$ws.Range("F3").Value = 0.2029300480133553
$ws.Range("F4").Value = 0.2248831547004854
$ws.Range("F5").Value = 0.1660892720754157
$ws.Range("F6").Value = 0.02706887246794999
$ws.Range("G2").Value = 0.4888330337345688
$ws.Range("G3").Value = 0.3772399484969685
$ws.Range("G4").Value = 0.2248831547004854
$ws.Range("G5").Value = 0.1916336854272317
$ws.Range("G6").Value = 0.05116366717998694
$ws.Range("H2").Value = 0.4888330337345688
$ws.Range("H3").Value = 0.4016239008754204
$ws.Range("H4").Value = 0.2513989312287693
$ws.Range("H5").Value = 0.3879584216252401
$ws.Range("H6").Value = 0.4096067714179726
$ws.Range("I2").Value = 0.511207746034443
$ws.Range("I3").Value = 0.4016239008754204
$ws.Range("I4").Value = 0.2987644177943506
$ws.Range("I5").Value = 0.3879584216252401
$ws.Range("I6").Value = 0.4096067714179726
$ws.Range("J2").Value = 0.511207746034443
$ws.Range("J3").Value = 0.4177439131883871
$ws.Range("J4").Value = 0.2987644177943506
$ws.Range("J5").Value = 0.3879584216252401
$ws.Range("J6").Value = 0.4163636449504306
$ws.Range("K2").Value = 0.511207746034443
$ws.Range("K3").Value = 0.4177439131883871
$ws.Range("K4").Value = 0.2991240157955214
$ws.Range("K5").Value = 0.3879584216252401
$ws.Range("K6").Value = 0.4163636449504306
$ws.Range("L2").Value = 0.5489807354121359
$ws.Range("L3").Value = 0.4177439131883871
$ws.Range("L4").Value = 0.4297434952983242
$ws.Range("L5").Value = 0.4334321809937014
$ws.Range("L6").Value = 0.4974639133422668
$ws.Range("M2").Value = 0.564723081965397
$ws.Range("M3").Value = 0.6789466782641659
$ws.Range("M4").Value = 0.6510064037951125
$ws.Range("M5").Value = 0.6378403740019559
$ws.Range("M6").Value = 0.6525642814121577
$ws.Range("N2").Value = 0.6031265739933629
$ws.Range("N3").Value = 0.6874892612813504
$ws.Range("N4").Value = 0.6510064037951125
$ws.Range("N5").Value = 0.6811961047480574
$ws.Range("N6").Value = 0.6837227027943986
$ws.Range("O2").Value = 0.7232631868768841
$ws.Range("O3").Value = 0.6934383488361134
$ws.Range("O4").Value = 0.6510064037951125
$ws.Range("O5").Value = 0.6811961047480574
$ws.Range("O6").Value = 0.6969355668759146
$ws.Range("P2").Value = 0.7232631868768841
$ws.Range("P3").Value = 0.7394799487808135
$ws.Range("P4").Value = 0.6510064037951125
$ws.Range("P5").Value = 0.6811961047480574
$ws.Range("P6").Value = 0.6969355668759146
$ws.Range("Q2").Value = 0.7232631868768841
$ws.Range("Q3").Value = 0.7416316321534959
$ws.Range("Q4").Value = 0.6510064037951125
$ws.Range("Q5").Value = 0.6914392663779447
$ws.Range("Q6").Value = 0.7145372056338682
$ws.Range("R2").Value = 0.728469278977491
$ws.Range("R3").Value = 0.7416316321534959
$ws.Range("R4").Value = 0.6510064037951125
$ws.Range("R5").Value = 0.6914392663779447
$ws.Range("R6").Value = 0.7145372056338682
$ws.Range("S2").Value = 0.728469278977491
$ws.Range("S3").Value = 0.7416316321534959
$ws.Range("S4").Value = 0.6510064037951125
$ws.Range("S5").Value = 0.6914392663779447
$ws.Range("S6").Value = 0.7145372056338682
$ws.Range("T2").Value = 0.728469278977491
$ws.Range("T3").Value = 0.7416316321534959
$ws.Range("T4").Value = 0.7675427910415907
$ws.Range("T5").Value = 0.7448224632134061
$ws.Range("T6").Value = 0.7835171813361947
$ws.Range("U2").Value = 0.8083243188805453
$ws.Range("U3").Value = 0.8164227276275569
$ws.Range("U4").Value = 0.855266227195162
$ws.Range("U5").Value = 0.7734584994668987
$ws.Range("U6").Value = 0.8195294526765603
$ws.Range("V2").Value = 0.8255261268332217
$ws.Range("V3").Value = 0.8328785857315624
$ws.Range("V4").Value = 0.855266227195162
$ws.Range("V5").Value = 0.8030280840166879
$ws.Range("V6").Value = 0.8507016354705971
$ws.Range("W2").Value = 0.8255261268332217
$ws.Range("W3").Value = 0.8489639678646116
$ws.Range("W4").Value = 0.855266227195162
$ws.Range("W5").Value = 0.8030280840166879
$ws.Range("W6").Value = 0.8507016354705971
$ws.Range("X2").Value = 0.8255261268332217
$ws.Range("X3").Value = 0.8489639678646116
$ws.Range("X4").Value = 0.871394506441258
$ws.Range("X5").Value = 0.8030280840166879
$ws.Range("X6").Value = 0.8507016354705971
$ws.Range("Y2").Value = 0.83391438309189
$ws.Range("Y3").Value = 0.883634244217353
$ws.Range("Y4").Value = 0.871394506441258
$ws.Range("Y5").Value = 0.8030280840166879
$ws.Range("Y6").Value = 0.8507016354705971
$ws.Range("Z2").Value = 0.83391438309189
$ws.Range("Z3").Value = 0.883634244217353
$ws.Range("Z4").Value = 0.8741626901624249
$ws.Range("Z5").Value = 0.8030280840166879
$ws.Range("Z6").Value = 0.8507016354705971

# --- Step3_DataPts_0.5 ---
$ws = $wb.Worksheets.Item("Step3_DataPts_0.5")
$ws.Range("C6").Value = 4
$ws.Range("D2").Value = 8
$ws.Range("D6").Value = 12
$ws.Range("E6").Value = 0.02706887246794999
$ws.Range("F2").Value = 0.511207746034443
$ws.Range("F3").Value = 0.6789466782641659
$ws.Range("F4").Value = 0.6510064037951125
$ws.Range("F5").Value = 0.6378403740019559
$ws.Range("F6").Value = 0.6525642814121577
$ws.Range("G2").Value = 6
$ws.Range("G6").Value = 8

# --- Step3_DataPts_0.7 ---
$ws = $wb.Worksheets.Item("Step3_DataPts_0.7")
$ws.Range("C6").Value = 4
$ws.Range("D3").Value = 15
$ws.Range("D4").Value = 19
$ws.Range("D5").Value = 19
$ws.Range("D6").Value = 16
$ws.Range("E6").Value = 0.02706887246794999
$ws.Range("F2").Value = 0.7232631868768841
$ws.Range("F3").Value = 0.7394799487808135
$ws.Range("F4").Value = 0.7675427910415907
$ws.Range("F5").Value = 0.7448224632134061
$ws.Range("F6").Value = 0.7145372056338682
$ws.Range("G3").Value = 12
$ws.Range("G4").Value = 17
$ws.Range("G5").Value = 17
$ws.Range("G6").Value = 12

# --- Step3_DataPts_0.8 ---
$ws = $wb.Worksheets.Item("Step3_DataPts_0.8")
$ws.Range("C6").Value = 4
$ws.Range("D2").Value = 20
$ws.Range("D3").Value = 20
$ws.Range("D4").Value = 20
$ws.Range("D5").Value = 21
$ws.Range("D6").Value = 20
$ws.Range("E6").Value = 0.02706887246794999
$ws.Range("F2").Value = 0.8083243188805453
$ws.Range("F3").Value = 0.8164227276275569
$ws.Range("F4").Value = 0.855266227195162
$ws.Range("F5").Value = 0.8030280840166879
$ws.Range("F6").Value = 0.8195294526765603
$ws.Range("G2").Value = 18
$ws.Range("G3").Value = 17
$ws.Range("G4").Value = 18
$ws.Range("G5").Value = 19
$ws.Range("G6").Value = 16

# --- Step3_DataPts_0.9 ---
$ws = $wb.Worksheets.Item("Step3_DataPts_0.9")
$ws.Range("C6").Value = 4
$ws.Range("D3").Value = 27
$ws.Range("D4").Value = 29
$ws.Range("D6").Value = 30
$ws.Range("E6").Value = 0.02706887246794999
$ws.Range("F2").Value = 0.9501377980429535
$ws.Range("F3").Value = 0.9015835409192469
$ws.Range("F4").Value = 0.9305017185503967
$ws.Range("F5").Value = 0.9237285577042523
$ws.Range("F6").Value = 0.9095014294734129
$ws.Range("G3").Value = 24
$ws.Range("G4").Value = 27
$ws.Range("G6").Value = 26
